# Refreshing content for 2019 revisions
# Applies the Parts_Actual_Costs.xlsx update:
#  - removes the "USB Micro-B Power Breakout" line item (old row 8), leaving row 8 blank
#  - replaces the "Waterproof airtight survival case" / "5v 5600 mAh USB battery" line
#    items (old rows 13/14) with new vendors/costs (3.3V 18650 batteries / PVC caps &
#    hose clamps)
#  - adds a new "Pololu voltage regulator" line item in row 16
#  - adds a new column J "Per Unit Shipping" = Shipping+tax / Quantity for every line item

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: delete the "USB Micro-B Power Breakout" entry (leave row blank) ---
$ws.Range("A8:I8").ClearContents()

# --- Header row -----------------------------------------------------------
$ws.Range("J1").Value = "Per Unit Shipping"

# --- Row 16: new "Pololu voltage regulator" line item -----------------------
$ws.Range("A16").Value = "Pololu voltage regulator"
$ws.Range("C16").Value = "Pololu"
$ws.Range("D16").Value = 87.25
$ws.Range("E16").Value = 12.37
$ws.Range("F16").Value = 25

# --- Row 14: was "5v 5600 mAh USB battery" / Amazon ------------------------
#            now "3.3V 18650 3500mAh batteries" / Imr Batteries
$ws.Range("A14").Value = "3.3V 18650 3500mAh batteries"
$ws.Range("B14").ClearContents()
$ws.Range("C14").Value = "Imr Batteries"
$ws.Range("D14").Value = 187.5
$ws.Range("E14").Value = 14.78
$ws.Range("F14").Value = 30

# --- Row 13: was "Waterproof airtight survival case" / Amazon -------------
#            now "PVC Caps and hose clamps" / Supplyhouse.com
$ws.Range("C13").Value = "Supplyhouse.com"
$ws.Range("A13").Value = "PVC Caps and hose clamps"
$ws.Range("B13").ClearContents()
$ws.Range("D13").Value = 71.45
$ws.Range("E13").Value = 38.7
$ws.Range("F13").Value = 25

# --- Column J: Per Unit Shipping = Shipping+tax / Quantity, rows 2-16 -----
$ws.Range("J2").Formula = "=E2/F2"
$ws.Range("J3").Formula = "=E3/F3"
$ws.Range("J4").Formula = "=E4/F4"
$ws.Range("J5").Formula = "=E5/F5"
$ws.Range("J6").Formula = "=E6/F6"
$ws.Range("J7").Formula = "=E7/F7"
$ws.Range("J9").Formula = "=E9/F9"
$ws.Range("J10").Formula = "=E10/F10"
$ws.Range("J11").Formula = "=E11/F11"
$ws.Range("J12").Formula = "=E12/F12"
$ws.Range("J13").Formula = "=E13/F13"
$ws.Range("J14").Formula = "=E14/F14"
$ws.Range("J15").Formula = "=E15/F15"
$ws.Range("J16").Formula = "=E16/F16"

$ws.Calculate()
